$d = $word.ActiveDocument

# 1) Insert a new paragraph right after the first heading
#    ("Windows Microsoft Excel Shortcuts") containing the new sentence.
#    Inserting *before* the existing (empty, unformatted) second paragraph
#    keeps the new paragraph free of the heading's character formatting.
$secondPara = $d.Paragraphs.Item(2)
$secondPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs.Item(2)
$newPara.Range.Text = "These are good key shorts that help a ton with financial modeling."

# 2) Color both section headings red (FF0000), keeping their existing
#    size/underline formatting. Apply via the whole paragraph Range (which
#    includes the paragraph mark) so both the run and the paragraph-mark
#    run-properties pick up the color, matching how Word colors a
#    fully-selected heading line.
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd("`r")
    if (($t -eq "Windows Microsoft Excel Shortcuts") -or ($t -eq "Windows Microsoft Excel Formulas")) {
        $p.Range.Font.Color = 255
    }
}

# 3) Collapse the "AltD + T" / "AltD + GG/GU" / comparison-operator runs
#    (which were split by spelling/grammar proofErr markers) back into a
#    single plain run each by replacing the text with itself.
$d.Content.Find.Execute(" AltD + T", $true, $false, $false, $false, $false, $true, 1, $false, " AltD + T", 2)
$d.Content.Find.Execute(" AltD + GG/GU", $true, $false, $false, $false, $false, $true, 1, $false, " AltD + GG/GU", 2)
$d.Content.Find.Execute("Comparison operators: <=  ,  >=  ,  =", $true, $false, $false, $false, $false, $true, 1, $false, "Comparison operators: <=  ,  >=  ,  =", 2)
